# Insert a new data row at row 172 (pushing the existing rows 172:237 down
# to 173:238) and populate it with a new "Zapallo italiano" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("172:172").Insert()

$ws.Range("A172").Value = 10
$ws.Range("B172").Value = "Vega Modelo de Temuco"
$ws.Range("C172").Value = "La Araucanía"
$ws.Range("D172").Value = 44468
$ws.Range("E172").Value = 9
$ws.Range("F172").Value = 100112032
$ws.Range("G172").Value = "Zapallo italiano"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 140
$ws.Range("K172").Value = 18000
$ws.Range("L172").Value = 20000
$ws.Range("M172").Value = 18714
$ws.Range("N172").Value = "`$/caja 60 unidades"
$ws.Range("O172").Value = "Región de Arica y Parinacota"
$ws.Range("P172").Value = 312
$ws.Range("Q172").Value = 60
$ws.Range("R172").Value = "Hortaliza"
